$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header for new column I ---
$ws.Range("I1").Value = "filter to get min and max bounds for Eet"

# --- 2. Copy F2:F72 (ann / ave quarters) values into I2:I72, then sort ascending ---
$vals = @()
for ($r = 2; $r -le 72; $r++) {
    $vals += $ws.Cells.Item($r, 6).Value2
}
$sortedVals = $vals | Sort-Object
for ($i = 0; $i -lt $sortedVals.Count; $i++) {
    $ws.Cells.Item($i + 2, 9).Value2 = $sortedVals[$i]
}

# --- 3. Turn on filtering for the new column and sort it ascending (keeps filter state) ---
$filterRange = $ws.Range("I1:I72")
$filterRange.AutoFilter()
$sortTarget = $ws.Range("I2:I72")
$sortTarget.Sort($ws.Range("I1"), 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 1)

# Register the hidden _FilterDatabase defined name that Excel creates for the filtered range
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$I`$1:`$I`$72")
$fdName.Visible = $false

# --- 4. Label + summary statistics for the 1990-2019 subset (rows 2-31 of column I) ---
$ws.Range("H74").Value = "From 1990 to 2019"

$ws.Range("H75").Value = "median"
$ws.Range("I75").Formula = "=MEDIAN(`$I`$2:`$I`$31)"

$ws.Range("H76").Value = "average"
$ws.Range("I76").Formula = "=AVERAGE(`$I`$2:`$I`$31)"

$ws.Range("H77").Value = "min"
$ws.Range("I77").Formula = "=MIN(`$I`$2:`$I`$31)"

$ws.Range("H78").Value = "max"
$ws.Range("I78").Formula = "=MAX(`$I`$2:`$I`$31)"

# Apply the Percent cell style to the new summary values
$ws.Range("I75:I78").Style = "Percent"

# --- 5. Sheet view state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 1
$ws.Range("I79").Select()
